# DataSource - Emision Motor - Prendario.xlsx
# Update test-run data (environment endpoint, credentials, policy year,
# pledge/loan dates, loan numbers) and drop the now-unused external
# workbook link + its hyperlink, per the commit:
# "se modificaron Datos y de PC_gestionDocumental se creo el test run
#  para emitir Motor Compl Blank Prendario y VariosAutos"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ambiente / credenciales (A2:D2) ---
# (A2 carries the quotePrefix style, so re-assert it as literal text)
$ws.Range("A2").Value = "'i-preproducciongestion.segurossura.com.ar"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("D2").Value = "silverarrow"

# --- NroCuenta (E2) ---
$ws.Range("E2").Value = 2033795924

# --- FechaInicio (Q2) / Anio (R2) ---
$ws.Range("Q2").Value = "'08/04/2021"
$ws.Range("R2").Value = 2021

# --- InicioDelPrestamo (AC2) ---
$ws.Range("AC2").Value = "'08/04/2021"

# --- Datos del prendario: columnas X,Y ahora son literales, no fórmulas
#     contra el libro externo [1]Emisión_Motor ---
$ws.Range("X2").Value = "1234567RGA012"
$ws.Range("Y2").Value = "1234567RGA012"

# --- FinDelPrestamo / Vencimiento1Cuota (AD2:AE2) ---
$ws.Range("AD2").Value = "'30/05/2021"
$ws.Range("AE2").Value = "'20/05/2021"

# --- NumPrestamo literal en W2 (antes fórmula contra el libro externo) ---
$ws.Range("W2").Value = "RGA012"

# El hipervínculo de B2 apuntaba al ambiente viejo; ya no corresponde
$ws.Hyperlinks.Delete()

# Ya no se usa el libro externo MaestroAltaPolizas.xlsx (W2:Y2 quedaron
# como literales), así que se rompe/quita el vínculo externo
$wb.BreakLink("MaestroAltaPolizas.xlsx", 1)

# Selección activa tal cual quedó en el archivo final
$null = $ws.Range("W7").Select()

$wb.Save()
